$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.768214821815491
$ws.Range("B1").Value = 1.956848740577698
$ws.Range("C1").Value = 2.32818865776062
$ws.Range("D1").Value = 2.889379501342773
$ws.Range("E1").Value = 3.673007249832153
